$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 01:32"

# Update rows whose data changed due to the daily refresh (new case counts,
# and re-sorting by "Casos totales" descending causes some countries to swap rows)

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 5093634
$ws.Cells.Item(4, 3).Value = 61356
$ws.Cells.Item(4, 4).Value = 2615683
$ws.Cells.Item(4, 5).Value = 2313905
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 1242
$ws.Cells.Item(4, 8).Value = 164046

# Row 5: Brasil
$ws.Cells.Item(5, 1).Value = "Brasil"
$ws.Cells.Item(5, 2).Value = 2967064
$ws.Cells.Item(5, 3).Value = 49502
$ws.Cells.Item(5, 4).Value = 2068394
$ws.Cells.Item(5, 5).Value = 798968
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 1058
$ws.Cells.Item(5, 8).Value = 99702

# Row 12: Colombia
$ws.Cells.Item(12, 1).Value = "Colombia"
$ws.Cells.Item(12, 2).Value = 367196
$ws.Cells.Item(12, 3).Value = 9486
$ws.Cells.Item(12, 4).Value = 198495
$ws.Cells.Item(12, 5).Value = 156451
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 311
$ws.Cells.Item(12, 8).Value = 12250

# Row 13: España
$ws.Cells.Item(13, 1).Value = "España"
$ws.Cells.Item(13, 2).Value = 361442
$ws.Cells.Item(13, 3).Value = 4507
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 3
$ws.Cells.Item(13, 8).Value = 28503

# Row 27: Canada
$ws.Cells.Item(27, 1).Value = "Canada"
$ws.Cells.Item(27, 2).Value = 118984
$ws.Cells.Item(27, 3).Value = 423
$ws.Cells.Item(27, 4).Value = 103433
$ws.Cells.Item(27, 5).Value = 6581
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 4
$ws.Cells.Item(27, 8).Value = 8970

# Row 39: Panama
$ws.Cells.Item(39, 1).Value = "Panama"
$ws.Cells.Item(39, 2).Value = 72560
$ws.Cells.Item(39, 3).Value = 1142
$ws.Cells.Item(39, 4).Value = 46675
$ws.Cells.Item(39, 5).Value = 24294
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 17
$ws.Cells.Item(39, 8).Value = 1591

# Row 40: Belgica
$ws.Cells.Item(40, 1).Value = "Belgica"
$ws.Cells.Item(40, 2).Value = 72016
$ws.Cells.Item(40, 3).Value = 858
$ws.Cells.Item(40, 4).Value = 17700
$ws.Cells.Item(40, 5).Value = 44455
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 2
$ws.Cells.Item(40, 8).Value = 9861

# Row 51: Nigeria
$ws.Cells.Item(51, 1).Value = "Nigeria"
$ws.Cells.Item(51, 2).Value = 45687
$ws.Cells.Item(51, 3).Value = 443
$ws.Cells.Item(51, 4).Value = 32637
$ws.Cells.Item(51, 5).Value = 12114
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 6
$ws.Cells.Item(51, 8).Value = 936

# Row 52: Japon
$ws.Cells.Item(52, 1).Value = "Japon"
$ws.Cells.Item(52, 2).Value = 43815
$ws.Cells.Item(52, 3).Value = 1552
$ws.Cells.Item(52, 4).Value = 30153
$ws.Cells.Item(52, 5).Value = 12629
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 7
$ws.Cells.Item(52, 8).Value = 1033

# Row 53: Barein
$ws.Cells.Item(53, 1).Value = "Barein"
$ws.Cells.Item(53, 2).Value = 43307
$ws.Cells.Item(53, 3).Value = 418
$ws.Cells.Item(53, 4).Value = 40276
$ws.Cells.Item(53, 5).Value = 2872
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 3
$ws.Cells.Item(53, 8).Value = 159

# Row 67: Venezuela
$ws.Cells.Item(67, 1).Value = "Venezuela"
$ws.Cells.Item(67, 2).Value = 24166
$ws.Cells.Item(67, 3).Value = 886
$ws.Cells.Item(67, 4).Value = 12470
$ws.Cells.Item(67, 5).Value = 11488
$ws.Cells.Item(67, 6).Value = 0
$ws.Cells.Item(67, 7).Value = 6
$ws.Cells.Item(67, 8).Value = 208

# Row 74: Chequia
$ws.Cells.Item(74, 1).Value = "Chequia"
$ws.Cells.Item(74, 2).Value = 18060
$ws.Cells.Item(74, 3).Value = 329
$ws.Cells.Item(74, 4).Value = 12749
$ws.Cells.Item(74, 5).Value = 4922
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 389

# Row 83: Sudan
$ws.Cells.Item(83, 1).Value = "Sudan"
$ws.Cells.Item(83, 2).Value = 11894
$ws.Cells.Item(83, 3).Value = 114
$ws.Cells.Item(83, 4).Value = 6243
$ws.Cells.Item(83, 5).Value = 4878
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 10
$ws.Cells.Item(83, 8).Value = 773

# Row 86: Noruega
$ws.Cells.Item(86, 1).Value = "Noruega"
$ws.Cells.Item(86, 2).Value = 9551
$ws.Cells.Item(86, 3).Value = 83
$ws.Cells.Item(86, 4).Value = 8857
$ws.Cells.Item(86, 5).Value = 438
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 256

# Row 89: Guayana Francesa
$ws.Cells.Item(89, 1).Value = "Guayana Francesa"
$ws.Cells.Item(89, 2).Value = 8204
$ws.Cells.Item(89, 3).Value = 77
$ws.Cells.Item(89, 4).Value = 7320
$ws.Cells.Item(89, 5).Value = 837
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 47

# Row 97: Mauritania
$ws.Cells.Item(97, 1).Value = "Mauritania"
$ws.Cells.Item(97, 2).Value = 6498
$ws.Cells.Item(97, 3).Value = 25
$ws.Cells.Item(97, 4).Value = 5443
$ws.Cells.Item(97, 5).Value = 898
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 157

# Row 113: Congo
$ws.Cells.Item(113, 1).Value = "Congo"
$ws.Cells.Item(113, 2).Value = 3637
$ws.Cells.Item(113, 3).Value = 91
$ws.Cells.Item(113, 4).Value = 1589
$ws.Cells.Item(113, 5).Value = 1990
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 58

# Row 114: Montenegro
$ws.Cells.Item(114, 1).Value = "Montenegro"
$ws.Cells.Item(114, 2).Value = 3549
$ws.Cells.Item(114, 3).Value = 69
$ws.Cells.Item(114, 4).Value = 2296
$ws.Cells.Item(114, 5).Value = 1192
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 1
$ws.Cells.Item(114, 8).Value = 61

# Row 117: Mayotte
$ws.Cells.Item(117, 1).Value = "Mayotte"
$ws.Cells.Item(117, 2).Value = 3068
$ws.Cells.Item(117, 3).Value = 26
$ws.Cells.Item(117, 4).Value = 2835
$ws.Cells.Item(117, 5).Value = 194
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 39

# Row 118: Suazilandia
$ws.Cells.Item(118, 1).Value = "Suazilandia"
$ws.Cells.Item(118, 2).Value = 3036
$ws.Cells.Item(118, 3).Value = 68
$ws.Cells.Item(118, 4).Value = 1476
$ws.Cells.Item(118, 5).Value = 1504
$ws.Cells.Item(118, 6).Value = 0
$ws.Cells.Item(118, 7).Value = 1
$ws.Cells.Item(118, 8).Value = 56

# Row 140: Uruguay
$ws.Cells.Item(140, 1).Value = "Uruguay"
$ws.Cells.Item(140, 2).Value = 1325
$ws.Cells.Item(140, 3).Value = 7
$ws.Cells.Item(140, 4).Value = 1095
$ws.Cells.Item(140, 5).Value = 193
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 37

# Row 156: Santo Tome y Principe
$ws.Cells.Item(156, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(156, 2).Value = 878
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(156, 4).Value = 799
$ws.Cells.Item(156, 5).Value = 64
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 15

# Row 157: Bahamas
$ws.Cells.Item(157, 1).Value = "Bahamas"
$ws.Cells.Item(157, 2).Value = 830
$ws.Cells.Item(157, 3).Value = 69
$ws.Cells.Item(157, 4).Value = 95
$ws.Cells.Item(157, 5).Value = 721
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 14

# Row 158: Botsuana
$ws.Cells.Item(158, 1).Value = "Botsuana"
$ws.Cells.Item(158, 2).Value = 804
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 63
$ws.Cells.Item(158, 5).Value = 739
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 2

# Row 159: Vietnam
$ws.Cells.Item(159, 1).Value = "Vietnam"
$ws.Cells.Item(159, 2).Value = 789
$ws.Cells.Item(159, 3).Value = 42
$ws.Cells.Item(159, 4).Value = 395
$ws.Cells.Item(159, 5).Value = 384
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 10

# Row 163: Reunion
$ws.Cells.Item(163, 1).Value = "Reunion"
$ws.Cells.Item(163, 2).Value = 675
$ws.Cells.Item(163, 3).Value = 4
$ws.Cells.Item(163, 4).Value = 631
$ws.Cells.Item(163, 5).Value = 39
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 5

# Row 175: Guadalupe
$ws.Cells.Item(175, 1).Value = "Guadalupe"
$ws.Cells.Item(175, 2).Value = 290
$ws.Cells.Item(175, 3).Value = 11
$ws.Cells.Item(175, 4).Value = 186
$ws.Cells.Item(175, 5).Value = 90
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 14

# Row 176: Eritrea
$ws.Cells.Item(176, 1).Value = "Eritrea"
$ws.Cells.Item(176, 2).Value = 285
$ws.Cells.Item(176, 3).Value = 3
$ws.Cells.Item(176, 4).Value = 245
$ws.Cells.Item(176, 5).Value = 40
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 0

# Row 187: Barbados
$ws.Cells.Item(187, 1).Value = "Barbados"
$ws.Cells.Item(187, 2).Value = 138
$ws.Cells.Item(187, 3).Value = 5
$ws.Cells.Item(187, 4).Value = 100
$ws.Cells.Item(187, 5).Value = 31
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 7

# Row 202: Timor Oriental
$ws.Cells.Item(202, 1).Value = "Timor Oriental"
$ws.Cells.Item(202, 2).Value = 25
$ws.Cells.Item(202, 3).Value = 0
$ws.Cells.Item(202, 4).Value = 24
$ws.Cells.Item(202, 5).Value = 1
$ws.Cells.Item(202, 6).Value = 0
$ws.Cells.Item(202, 7).Value = 0
$ws.Cells.Item(202, 8).Value = 0

# Row 203: Santa Lucia
$ws.Cells.Item(203, 1).Value = "Santa Lucia"
$ws.Cells.Item(203, 2).Value = 25
$ws.Cells.Item(203, 3).Value = 0
$ws.Cells.Item(203, 4).Value = 24
$ws.Cells.Item(203, 5).Value = 1
$ws.Cells.Item(203, 6).Value = 0
$ws.Cells.Item(203, 7).Value = 0
$ws.Cells.Item(203, 8).Value = 0
